$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the "3: void loop..." paragraph to the
#    end of the "Level 2: Using External Documentation" paragraph.
# ---------------------------------------------------------------------------

# Remove the old bookmark (this only deletes the bookmark markers, not text).
$oldMark = $d.Bookmarks("_GoBack")
$oldMark.Delete()

# Locate the "Level 2: Using External Documentation" paragraph and collapse a
# range to right after its text (i.e. just before the paragraph mark). A
# range collapsed exactly on a paragraph boundary confuses Bookmarks.Add, so
# we temporarily append a marker character, bookmark just before it, then
# remove the marker again; the bookmark sticks at the correct spot.
$target = $d.Content
$target.Find.Execute("Level 2: Using External Documentation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target.Collapse(0)
$target.InsertAfter("@@MARKER@@")

$newSpot = $d.Content
$newSpot.Find.Execute("@@MARKER@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$newSpot.Collapse(1)

$d.Bookmarks.Add("_GoBack", $newSpot)

$markerRange = $d.Content
$markerRange.Find.Execute("@@MARKER@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$markerRange.Delete()

# ---------------------------------------------------------------------------
# 2) Remove the trailing duplicate "void loop / digital Write / delay" block
#    near the end of the document, keeping a single blank paragraph in its
#    place (right before the final section properties).
# ---------------------------------------------------------------------------

# The duplicate block is the *last* occurrence of "void loop" in the body
# (the first occurrence is inside the "3: void loop..." paragraph earlier in
# the document, which must be left untouched). Walk paragraphs from the end
# to find it reliably instead of relying on Find (which matches first hit).
$startParaIndex = -1
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    if ($d.Paragraphs($i).Range.Text -like "void loop*") {
        $startParaIndex = $i
        break
    }
}

$keepParaIndex = $startParaIndex + 3
$pKeepBefore = $d.Paragraphs($keepParaIndex)

$deleteRange1 = $d.Range($d.Paragraphs($startParaIndex).Range.Start, $pKeepBefore.Range.Start)
$deleteRange1.Delete()

# Re-resolve the kept paragraph *after* the structural edit above rather than
# reusing the pre-edit Paragraph object, which does not track the original
# content once earlier paragraphs are removed.
$pKeepAfter = $d.Paragraphs($startParaIndex)
$deleteRange2 = $d.Range($pKeepAfter.Range.End, $d.Content.End)
$deleteRange2.Delete()
